$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 3017
$ws.Range("I21").Value = 3017
$ws.Range("K21").Value = 3017
$ws.Range("M21").Value = -2549

$ws.Range("H23").Value = 3017
$ws.Range("I23").Value = 3017
$ws.Range("K23").Value = 3017
$ws.Range("M23").Value = -2783

$ws.Range("H116").Value = 4259
$ws.Range("I116").Value = 3598.3333
$ws.Range("K116").Value = 3598.3333
$ws.Range("M116").Value = -156.3332999999998

$ws.Range("H125").Value = 127653.375
$ws.Range("I125").Value = 1031
$ws.Range("J125").Value = 203626.8
$ws.Range("K125").Value = 9279
$ws.Range("L125").Value = 1832641.2
$ws.Range("M125").Value = -6819
$ws.Range("N125").Value = -1837561.2


# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 895
$ws.Range("I2").Value = 895
$ws.Range("K2").Value = 895
$ws.Range("M2").Value = -782

$ws.Range("H37").Value = 9999.571
$ws.Range("J37").Value = 9999.571
$ws.Range("L37").Value = 9999.571
$ws.Range("N37").Value = -10545.571

$ws.Range("H88").Value = 1587.4762
$ws.Range("I88").Value = 394.5
$ws.Range("K88").Value = 394.5
$ws.Range("M88").Value = 11.5

$ws.Range("H91").Value = 1587.4762
$ws.Range("I91").Value = 394.5
$ws.Range("K91").Value = 394.5
$ws.Range("M91").Value = 1009.5

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws.Range("H116").Value = 895
$ws.Range("I116").Value = 895
$ws.Range("K116").Value = 895
$ws.Range("M116").Value = 1399

$ws.Range("H122").Value = 2893.2
$ws.Range("I122").Value = 2893.2
$ws.Range("K122").Value = 8679.599999999999
$ws.Range("M122").Value = -6229.599999999999


# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 895
$ws.Range("I3").Value = 895
$ws.Range("K3").Value = 895
$ws.Range("M3").Value = -781


# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1832
$ws.Range("I16").Value = 1499.25
$ws.Range("K16").Value = 1499.25
$ws.Range("M16").Value = -1212.25

$ws.Range("H55").Value = 11550.5
$ws.Range("J55").Value = 11550.5
$ws.Range("L55").Value = 11550.5
$ws.Range("N55").Value = -12180.5

$ws.Range("H58").Value = 3883.2222
$ws.Range("I58").Value = 2866.5
$ws.Range("K58").Value = 2866.5
$ws.Range("M58").Value = -2663.5

$ws.Range("H92").Value = 80000
$ws.Range("J92").Value = 80000
$ws.Range("L92").Value = 80000
$ws.Range("N92").Value = -84992

$ws.Range("H99").Value = 4490.2
$ws.Range("I99").Value = 3979.6
$ws.Range("K99").Value = 3979.6
$ws.Range("M99").Value = -2481.6

$ws.Range("H113").Value = 1832
$ws.Range("I113").Value = 1499.25
$ws.Range("K113").Value = 1499.25
$ws.Range("M113").Value = 670.75

$ws.Range("H122").Value = 2626
$ws.Range("I122").Value = 2626
$ws.Range("K122").Value = 7878
$ws.Range("M122").Value = -5428

$ws.Range("H126").Value = 4490.2
$ws.Range("I126").Value = 3979.6
$ws.Range("K126").Value = 11938.8
$ws.Range("M126").Value = -9468.799999999999

$ws.Range("H134").Value = 1161.7142
$ws.Range("I134").Value = 1161.7142
$ws.Range("K134").Value = 3485.1426
$ws.Range("M134").Value = -950.1425999999997

$ws.Range("H136").Value = 3883.2222
$ws.Range("I136").Value = 2866.5
$ws.Range("K136").Value = 8599.5
$ws.Range("M136").Value = -6049.5


# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H132").Value = 3431.25
$ws.Range("J132").Value = 3364.2856
$ws.Range("L132").Value = 30278.5704
$ws.Range("N132").Value = -35338.5704


# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 979
$ws.Range("I113").Value = 979
$ws.Range("K113").Value = 979
$ws.Range("M113").Value = 1191

$ws.Range("H122").Value = 1014.875
$ws.Range("I122").Value = 1002.5
$ws.Range("J122").Value = 1027.25
$ws.Range("K122").Value = 3007.5
$ws.Range("L122").Value = 3081.75
$ws.Range("M122").Value = -557.5
$ws.Range("N122").Value = -7981.75

$ws.Range("H126").Value = 13759.8
$ws.Range("I126").Value = 9600
$ws.Range("K126").Value = 28800
$ws.Range("M126").Value = -26330


# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11695.77
$ws.Range("I7").Value = 11267.728
$ws.Range("K7").Value = 11267.728
$ws.Range("M7").Value = -11155.728

$ws.Range("H40").Value = 5675.643
$ws.Range("I40").Value = 5359.909
$ws.Range("K40").Value = 5359.909
$ws.Range("M40").Value = -5223.909

$ws.Range("H46").Value = 2387.6875
$ws.Range("I46").Value = 2595.875
$ws.Range("J46").Value = 2179.5
$ws.Range("K46").Value = 2595.875
$ws.Range("L46").Value = 2179.5
$ws.Range("M46").Value = -2407.875
$ws.Range("N46").Value = -2555.5

$ws.Range("H55").Value = 235.33333
$ws.Range("I55").Value = 150
$ws.Range("J55").Value = 252.4
$ws.Range("K55").Value = 150
$ws.Range("L55").Value = 252.4
$ws.Range("M55").Value = 23
$ws.Range("N55").Value = -598.4

$ws.Range("H61").Value = 2131.1667
$ws.Range("I61").Value = 2317.4
$ws.Range("K61").Value = 2317.4
$ws.Range("M61").Value = -2115.4

$ws.Range("H113").Value = 2131.1667
$ws.Range("I113").Value = 2317.4
$ws.Range("K113").Value = 2317.4
$ws.Range("M113").Value = -147.4000000000001

$ws.Range("H122").Value = 3401.75
$ws.Range("I122").Value = 3244.8572
$ws.Range("K122").Value = 9734.571599999999
$ws.Range("M122").Value = -7284.571599999999

$ws.Range("H126").Value = 11695.77
$ws.Range("I126").Value = 11267.728
$ws.Range("K126").Value = 33803.18399999999
$ws.Range("M126").Value = -31333.18399999999


# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 85999.5
$ws.Range("I126").Value = 72000
$ws.Range("K126").Value = 216000
$ws.Range("M126").Value = -213530

$ws.Range("H132").Value = 1761.4546
$ws.Range("I132").Value = 1761.4546
$ws.Range("K132").Value = 5284.3638
$ws.Range("M132").Value = -2754.3638

$ws.Range("H136").Value = 4354.8
$ws.Range("I136").Value = 6499.8335
$ws.Range("J136").Value = 1137.25
$ws.Range("K136").Value = 19499.5005
$ws.Range("L136").Value = 3411.75
$ws.Range("M136").Value = -16949.5005
$ws.Range("N136").Value = -8511.75

